# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计" (i.e. as the 2nd tab),
#   pushing 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q1 / 2020-Q4 one slot later.
# - Populate "2022-Q3" with the single fund row for 中金丰硕混合 (005396).
# - Update the "总计" summary sheet: insert a new row for 2022-Q3
#   (1 fund, 0.1 billion) above the existing rows, shifting the rest down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet in place.
#    Final target layout (rows 2..7), all columns already carry the right
#    types/styles on rows that pre-existed (A2:A6 use style "2"); only the
#    brand-new row 7 needs that style copied explicitly.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("A6").Copy() | Out-Null
$summary.Range("A7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$rows = @(
    @(0, "2022-Q3", 1, 0.1),
    @(1, "2022-Q2", 4, 0.24),
    @(2, "2022-Q1", 5, 0.53),
    @(3, "2021-Q4", 1, 3.36),
    @(4, "2021-Q1", 1, 0),
    @(5, "2020-Q4", 2, 1.08)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $summary.Cells.Item($r, 1).Value = $data[0]
    $summary.Cells.Item($r, 2).Value = $data[1]
    $summary.Cells.Item($r, 3).Value = $data[2]
    $summary.Cells.Item($r, 4).Value = $data[3]
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right before the current "2022-Q2"
#    tab (so it lands as the second sheet overall).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item("2022-Q2"))
$q3.Name = "2022-Q3"

# Re-resolve "2022-Q2" now that the sheet collection has shifted — a handle
# grabbed before Worksheets.Add() no longer tracks the live sheet.
$q2 = $wb.Worksheets.Item("2022-Q2")

# Borrow the header-row / first-data-row formatting from the sheet we just
# pushed down, so the new tab matches the house style (bold, bordered,
# centered headers + the styled A column) without hand-building a style.
$q2.Range("B1:H1").Copy() | Out-Null
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$q2.Range("A2").Copy() | Out-Null
$q3.Range("A2").PasteSpecial(-4122) | Out-Null         # xlPasteFormats

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0

# B2..G2 are stored as text in the source data (even the numeric-looking
# ones), so force text formatting before assigning to avoid Excel's
# automatic number coercion / loss of the leading zero in the fund code.
$q3.Range("B2:G2").NumberFormat = "@"
$q3.Range("B2").Value = "005396"
$q3.Range("C2").Value = "中金丰硕混合"
$q3.Range("D2").Value = "1.77"
$q3.Range("E2").Value = "76.61"
$q3.Range("F2").Value = "5.41"
$q3.Range("G2").Value = "0.0958"

$q3.Range("H2").Value = 6
